# ---------------------------------------------------------------------------
# B6-PowerPoint.pptx edit
#
# 1) Three tables (on the slides holding the BTEC "component 3" grading
#    grids) get re-pointed from the deck's bespoke table style
#    {50BFD30B-1374-4D6E-8D2F-C66EE08CED7B} to the built-in gallery style
#    {297E6AD9-8D99-41DD-8769-62CFE0B0D5B5}.
#
# 2) The presentation's theme color palette is switched from the custom
#    "Red Violet" palette (part of the "Integral" theme) to the stock
#    Office palette.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

$oldTableStyle = "{50BFD30B-1374-4D6E-8D2F-C66EE08CED7B}"
$newTableStyle = "{297E6AD9-8D99-41DD-8769-62CFE0B0D5B5}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldTableStyle) {
                $tbl.ApplyStyle($newTableStyle)
            }
        }
    }
}

# Re-color the theme from the "Red Violet" scheme to the stock "Office"
# scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink — in that order).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le $officeColors.Length; $k++) {
    $themeColors.Colors($k).RGB = $officeColors[$k - 1]
}
